$wb = $excel.ActiveWorkbook

# ---- Update existing "GLOBAL RESULTS" sheet (sheet1) ----
$ws1 = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws1.Range("A1").Value = 'Description'
$ws1.Range("B1").Value = 'Unit'
$ws1.Range("C1").Value = 'Value'
$ws1.Range("A2").Value = 'Reference Range'
$ws1.Range("B2").Value = 'nmi'
$ws1.Range("C2").Value = 825.0539956803455
$ws1.Range("A3").Value = 'Material density'
$ws1.Range("B3").Value = 'kg/m³'
$ws1.Range("C3").Value = 2711.0
$ws1.Range("A4").Value = 'Single passenger Mass'
$ws1.Range("B4").Value = 'kg'
$ws1.Range("C4").Value = 99.0
$ws1.Range("A5").Value = ' '
$ws1.Range("A6").Value = 'Maximum Take-Off Mass'
$ws1.Range("B6").Value = 'kg'
$ws1.Range("C6").Value = 25070.404491806272
$ws1.Range("A7").Value = 'Maximum Take-Off Weight'
$ws1.Range("B7").Value = 'N'
$ws1.Range("C7").Value = 245856.68220957194
$ws1.Range("A8").Value = 'Take-Off Mass'
$ws1.Range("B8").Value = 'kg'
$ws1.Range("C8").Value = 24690.447825139614
$ws1.Range("A9").Value = 'Take-Off Weight'
$ws1.Range("B9").Value = 'N'
$ws1.Range("C9").Value = 242130.58016440534
$ws1.Range("A10").Value = 'Maximum Landing Mass'
$ws1.Range("B10").Value = 'kg'
$ws1.Range("C10").Value = 22577.803042625652
$ws1.Range("A11").Value = 'Maximum Landing Weight'
$ws1.Range("B11").Value = 'N'
$ws1.Range("C11").Value = 221412.6122079648
$ws1.Range("A12").Value = 'Maximum Passengers Mass'
$ws1.Range("B12").Value = 'kg'
$ws1.Range("C12").Value = 7128.0
$ws1.Range("A13").Value = 'Maximum Passengers Weight'
$ws1.Range("B13").Value = 'N'
$ws1.Range("C13").Value = 69901.80119999999
$ws1.Range("A14").Value = 'Fuel Mass'
$ws1.Range("B14").Value = 'kg'
$ws1.Range("C14").Value = 4954.49177994087
$ws1.Range("A15").Value = 'Fuel Weight'
$ws1.Range("B15").Value = 'N'
$ws1.Range("C15").Value = 48586.96681375712
$ws1.Range("A16").Value = 'Crew Mass'
$ws1.Range("B16").Value = 'kg'
$ws1.Range("C16").Value = 229.54364550000003
$ws1.Range("A17").Value = 'Crew Weight'
$ws1.Range("B17").Value = 'N'
$ws1.Range("C17").Value = 2251.0541911425744
$ws1.Range("A18").Value = 'Maximum Zero Fuel Mass'
$ws1.Range("B18").Value = 'kg'
$ws1.Range("C18").Value = 20131.956045198745
$ws1.Range("A19").Value = 'Maximum Zero Fuel Weight'
$ws1.Range("B19").Value = 'N'
$ws1.Range("C19").Value = 197427.04675064824
$ws1.Range("A20").Value = 'Zero Fuel Mass'
$ws1.Range("B20").Value = 'kg'
$ws1.Range("C20").Value = 19735.956045198745
$ws1.Range("A21").Value = 'Zero Fuel Weight'
$ws1.Range("B21").Value = 'N'
$ws1.Range("C21").Value = 193543.61335064823
$ws1.Range("A22").Value = 'Operating Empty Mass'
$ws1.Range("B22").Value = 'kg'
$ws1.Range("C22").Value = 13003.956045198749
$ws1.Range("A23").Value = 'Operating Empty Weight'
$ws1.Range("B23").Value = 'N'
$ws1.Range("C23").Value = 127525.24555064828
$ws1.Range("A24").Value = 'Empty Mass'
$ws1.Range("B24").Value = 'kg'
$ws1.Range("C24").Value = 12774.412399698747
$ws1.Range("A25").Value = 'Empty Weight'
$ws1.Range("B25").Value = 'N'
$ws1.Range("C25").Value = 125274.19135950568
$ws1.Range("A26").Value = 'Manufacturer Empty Mass'
$ws1.Range("B26").Value = 'kg'
$ws1.Range("C26").Value = 12188.456399698749
$ws1.Range("A27").Value = 'Manufacturer Empty Weight'
$ws1.Range("B27").Value = 'N'
$ws1.Range("C27").Value = 119527.92595210572
$ws1.Range("A28").Value = 'Operating Item Mass'
$ws1.Range("B28").Value = 'kg'
$ws1.Range("C28").Value = 585.9559999999999
$ws1.Range("A29").Value = 'Operating Item Weight'
$ws1.Range("B29").Value = 'N'
$ws1.Range("C29").Value = 5746.265407399998
$ws1.Range("A30").Value = 'Trapped Fuel Oil Mass'
$ws1.Range("B30").Value = 'kg'
$ws1.Range("C30").Value = 0.0
$ws1.Range("A31").Value = 'Trapped Fuel Oil Weight'
$ws1.Range("B31").Value = 'N'
$ws1.Range("C31").Value = 0.0
$ws1.Range("A32").Value = 'Operating Empty Mass'
$ws1.Range("B32").Value = 'kg'
$ws1.Range("C32").Value = 13003.956045198749
$ws1.Range("A33").Value = 'Operating Empty Weight'
$ws1.Range("B33").Value = 'N'
$ws1.Range("C33").Value = 127525.24555064828

# Sheet with a fully-styled 4-column header row (A1:D1, style index 1)
# used as the format source for the new sheets below.
$wsHeaderSrc = $wb.Worksheets.Item("FUSELAGE")

# ---- Add new sheets in order: NACELLES, POWER PLANT, LANDING GEARS, SYSTEMS ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNacelles = $wb.Worksheets.Add($null, $lastSheet)
$wsNacelles.Name = "NACELLES"
$wsHeaderSrc.Range("A1:D1").Copy()
$wsNacelles.Range("A1:D1").PasteSpecial(-4122)
$wsNacelles.Range("A1").Value = 'Description'
$wsNacelles.Range("B1").Value = 'Unit'
$wsNacelles.Range("C1").Value = 'Value'
$wsNacelles.Range("D1").Value = 'Percent Error'
$wsNacelles.Range("A2").Value = 'Total Reference Mass'
$wsNacelles.Range("B2").Value = 'Kg'
$wsNacelles.Range("C2").Value = 818.7999999999997
$wsNacelles.Range("A3").Value = 'Total mass estimated'
$wsNacelles.Range("B3").Value = 'Kg'
$wsNacelles.Range("C3").Value = 645.9999999999998
$wsNacelles.Range("D3").Value = -21.10405471421592
$wsNacelles.Range("A4").Value = ' '
$wsNacelles.Range("A5").Value = 'WEIGHT ESTIMATION METHODS COMPARISON FOR EACH NACELLE'
$wsNacelles.Range("A6").Value = ' '
$wsNacelles.Range("A7").Value = 'NACELLE 1'
$wsNacelles.Range("A8").Value = 'Reference Mass'
$wsNacelles.Range("B8").Value = 'Kg'
$wsNacelles.Range("C8").Value = 409.4
$wsNacelles.Range("A9").Value = 'JENKINSON'
$wsNacelles.Range("B9").Value = 'Kg'
$wsNacelles.Range("C9").Value = 235.0
$wsNacelles.Range("D9").Value = -42.59892525647288
$wsNacelles.Range("A10").Value = 'TORENBEEK_1976'
$wsNacelles.Range("B10").Value = 'Kg'
$wsNacelles.Range("C10").Value = 540.0
$wsNacelles.Range("D10").Value = 31.90034196384954
$wsNacelles.Range("A11").Value = 'TORENBEEK_1982'
$wsNacelles.Range("B11").Value = 'Kg'
$wsNacelles.Range("C11").Value = 194.0
$wsNacelles.Range("D11").Value = -52.61358085002442
$wsNacelles.Range("A12").Value = 'Estimated Mass '
$wsNacelles.Range("B12").Value = 'Kg'
$wsNacelles.Range("C12").Value = 323.0
$wsNacelles.Range("A13").Value = ' '
$wsNacelles.Range("A14").Value = 'NACELLE 2'
$wsNacelles.Range("A15").Value = 'Reference Mass'
$wsNacelles.Range("B15").Value = 'Kg'
$wsNacelles.Range("C15").Value = 409.4
$wsNacelles.Range("A16").Value = 'JENKINSON'
$wsNacelles.Range("B16").Value = 'Kg'
$wsNacelles.Range("C16").Value = 235.0
$wsNacelles.Range("D16").Value = -42.59892525647288
$wsNacelles.Range("A17").Value = 'TORENBEEK_1976'
$wsNacelles.Range("B17").Value = 'Kg'
$wsNacelles.Range("C17").Value = 540.0
$wsNacelles.Range("D17").Value = 31.90034196384954
$wsNacelles.Range("A18").Value = 'TORENBEEK_1982'
$wsNacelles.Range("B18").Value = 'Kg'
$wsNacelles.Range("C18").Value = 194.0
$wsNacelles.Range("D18").Value = -52.61358085002442
$wsNacelles.Range("A19").Value = 'Estimated Mass '
$wsNacelles.Range("B19").Value = 'Kg'
$wsNacelles.Range("C19").Value = 323.0
$wsNacelles.Range("A20").Value = ' '

$wsPowerPlant = $wb.Worksheets.Add($null, $wsNacelles)
$wsPowerPlant.Name = "POWER PLANT"
$wsHeaderSrc.Range("A1:D1").Copy()
$wsPowerPlant.Range("A1:D1").PasteSpecial(-4122)
$wsPowerPlant.Range("A1").Value = 'Description'
$wsPowerPlant.Range("B1").Value = 'Unit'
$wsPowerPlant.Range("C1").Value = 'Value'
$wsPowerPlant.Range("D1").Value = 'Percent Error'
$wsPowerPlant.Range("A2").Value = 'Total Reference Mass'
$wsPowerPlant.Range("B2").Value = 'Kg'
$wsPowerPlant.Range("C2").Value = 965.2445633599998
$wsPowerPlant.Range("A3").Value = 'Total mass estimated'
$wsPowerPlant.Range("B3").Value = 'Kg'
$wsPowerPlant.Range("C3").Value = 1447.8668450399996
$wsPowerPlant.Range("D3").Value = -21.10405471421592
$wsPowerPlant.Range("A4").Value = ' '
$wsPowerPlant.Range("A5").Value = 'WEIGHT ESTIMATION METHODS COMPARISON FOR EACH ENGINE'
$wsPowerPlant.Range("A6").Value = ' '
$wsPowerPlant.Range("A7").Value = 'ENGINE 1'
$wsPowerPlant.Range("A8").Value = 'Reference Mass'
$wsPowerPlant.Range("B8").Value = 'Kg'
$wsPowerPlant.Range("C8").Value = 482.6222816799999
$wsPowerPlant.Range("A9").Value = 'Total Mass'
$wsPowerPlant.Range("B9").Value = 'Kg'
$wsPowerPlant.Range("C9").Value = 723.9334225199998
$wsPowerPlant.Range("A10").Value = ' '
$wsPowerPlant.Range("A11").Value = 'ENGINE 2'
$wsPowerPlant.Range("A12").Value = 'Reference Mass'
$wsPowerPlant.Range("B12").Value = 'Kg'
$wsPowerPlant.Range("C12").Value = 482.6222816799999
$wsPowerPlant.Range("A13").Value = 'Total Mass'
$wsPowerPlant.Range("B13").Value = 'Kg'
$wsPowerPlant.Range("C13").Value = 723.9334225199998
$wsPowerPlant.Range("A14").Value = ' '

$wsLandingGears = $wb.Worksheets.Add($null, $wsPowerPlant)
$wsLandingGears.Name = "LANDING GEARS"
$wsHeaderSrc.Range("A1:D1").Copy()
$wsLandingGears.Range("A1:D1").PasteSpecial(-4122)
$wsLandingGears.Range("A1").Value = 'Description'
$wsLandingGears.Range("B1").Value = 'Unit'
$wsLandingGears.Range("C1").Value = 'Value'
$wsLandingGears.Range("D1").Value = 'Percent Error'
$wsLandingGears.Range("A2").Value = 'Reference Mass'
$wsLandingGears.Range("B2").Value = 'Kg'
$wsLandingGears.Range("C2").Value = 675.8
$wsLandingGears.Range("A3").Value = 'Overall Mass'
$wsLandingGears.Range("B3").Value = 'Kg'
$wsLandingGears.Range("C3").Value = 987.8878901438338
$wsLandingGears.Range("A4").Value = ' '
$wsLandingGears.Range("A5").Value = 'WEIGHT ESTIMATION METHODS COMPARISON'
$wsLandingGears.Range("A6").Value = 'ROSKAM'
$wsLandingGears.Range("B6").Value = 'Kg'
$wsLandingGears.Range("C6").Value = 821.0
$wsLandingGears.Range("D6").Value = 21.48564664101806
$wsLandingGears.Range("A7").Value = 'STANFORD'
$wsLandingGears.Range("B7").Value = 'Kg'
$wsLandingGears.Range("C7").Value = 1003.0
$wsLandingGears.Range("D7").Value = 48.41669132879551
$wsLandingGears.Range("A8").Value = 'TORENBEEK_1982'
$wsLandingGears.Range("B8").Value = 'Kg'
$wsLandingGears.Range("C8").Value = 1135.0
$wsLandingGears.Range("D8").Value = 67.94909736608466
$wsLandingGears.Range("A9").Value = 'TORENBEEK_2013'
$wsLandingGears.Range("B9").Value = 'Kg'
$wsLandingGears.Range("C9").Value = 988.0
$wsLandingGears.Range("D9").Value = 46.19709973364902
$wsLandingGears.Range("A10").Value = 'Estimated Mass '
$wsLandingGears.Range("B10").Value = 'Kg'
$wsLandingGears.Range("C10").Value = 821.0

$wsSystems = $wb.Worksheets.Add($null, $wsLandingGears)
$wsSystems.Name = "SYSTEMS"
$wsHeaderSrc.Range("A1:D1").Copy()
$wsSystems.Range("A1:D1").PasteSpecial(-4122)
$wsSystems.Range("A1").Value = 'Description'
$wsSystems.Range("B1").Value = 'Unit'
$wsSystems.Range("C1").Value = 'Value'
$wsSystems.Range("D1").Value = 'Percent Error'
$wsSystems.Range("A2").Value = 'Reference Mass'
$wsSystems.Range("B2").Value = 'Kg'
$wsSystems.Range("C2").Value = 2118.0
$wsSystems.Range("A3").Value = 'Overall Mass'
$wsSystems.Range("B3").Value = 'Kg'
$wsSystems.Range("C3").Value = 2324.2562213254223
$wsSystems.Range("A4").Value = ' '
$wsSystems.Range("A5").Value = 'WEIGHT ESTIMATION METHODS COMPARISON'
$wsSystems.Range("A6").Value = 'TORENBEEK_2013'
$wsSystems.Range("B6").Value = 'Kg'
$wsSystems.Range("C6").Value = 2324.0
$wsSystems.Range("D6").Value = 9.726156751652502
$wsSystems.Range("A7").Value = 'Estimated Mass '
$wsSystems.Range("B7").Value = 'Kg'
$wsSystems.Range("C7").Value = 2324.0


